# Update the exam name for every student record to the new 2026 exam title.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F2:F5").Value = "TALENT SEARCH EXAMINATION 2026"

# All exam dates are now the same single date (01-Feb-2026) instead of a
# different year per row.
$ws.Range("H3").Value = "01-Feb-2026"
$ws.Range("H4").Value = "01-Feb-2026"
$ws.Range("H5").Value = "01-Feb-2026"

# Widen column F (examName) so the longer exam title is fully visible.
$ws.Columns.Item(6).ColumnWidth = 48

# Move the selection to G9 (matches the author's final cursor position).
$null = $ws.Range("G9").Select()
